$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AB2").Value = -981.2587840583742
$ws.Range("AC2").Value = "umolO2/min/m2"
$ws.Range("AD2").Value = -981.2587840583742
$ws.Range("D2").Value = 218.1463131112375
$ws.Range("E2").Value = -1.735163104968382
$ws.Range("F2").Value = 0.991
$ws.Range("H2").Value = 101
$ws.Range("I2").Value = 181
$ws.Range("J2").Value = 24.95
$ws.Range("K2").Value = 44.95
$ws.Range("L2").Value = 174.878177
$ws.Range("M2").Value = 138.524164
$ws.Range("N2").Value = -1.735163104968382
$ws.Range("O2").Value = -0.07817522212233576
$ws.Range("P2").Value = -1.656987882846046
$ws.Range("Q2").Value = -1.656987882846046
$ws.Range("T2").Value = 0.1450048780487805
$ws.Range("V2").Value = 0.0002448603057459146
$ws.Range("Z2").Value = -0.2402713258803979
$ws.Range("AB3").Value = -1666.799694145389
$ws.Range("AC3").Value = "umolO2/min/m2"
$ws.Range("AD3").Value = -1666.799694145389
$ws.Range("D3").Value = 237.7066000306082
$ws.Range("E3").Value = -1.741050521228543
$ws.Range("F3").Value = 0.99
$ws.Range("H3").Value = 101
$ws.Range("I3").Value = 181
$ws.Range("J3").Value = 24.98333333333333
$ws.Range("K3").Value = 44.98333333333333
$ws.Range("L3").Value = 193.530558
$ws.Range("M3").Value = 160.456873
$ws.Range("N3").Value = -1.741050521228543
$ws.Range("O3").Value = -0.07817522212233576
$ws.Range("P3").Value = -1.662875299106208
$ws.Range("Q3").Value = -1.662875299106208
$ws.Range("T3").Value = 0.1492487804878049
$ws.Range("V3").Value = 0.0001488973818309612
$ws.Range("Z3").Value = -0.2481821104948953
$ws.Range("AB4").Value = -1369.97470459495
$ws.Range("AC4").Value = "umolO2/min/m2"
$ws.Range("AD4").Value = -1369.97470459495
$ws.Range("D4").Value = 206.8458919471439
$ws.Range("E4").Value = -2.150765142818427
$ws.Range("F4").Value = 0.931
$ws.Range("H4").Value = 102
$ws.Range("I4").Value = 182
$ws.Range("J4").Value = 25.05
$ws.Range("K4").Value = 45.05
$ws.Range("L4").Value = 152.385555
$ws.Range("M4").Value = 113.876846
$ws.Range("N4").Value = -2.150765142818427
$ws.Range("O4").Value = -0.07817522212233576
$ws.Range("P4").Value = -2.072589920696091
$ws.Range("Q4").Value = -2.072589920696091
$ws.Range("T4").Value = 0.1469268292682927
$ws.Range("V4").Value = 0.0002222807942365138
$ws.Range("Z4").Value = -0.3045190654212989
$ws.Range("AB5").Value = -1034.295753429247
$ws.Range("AC5").Value = "umolO2/min/m2"
$ws.Range("AD5").Value = -1034.295753429247
$ws.Range("D5").Value = 259.642197738347
$ws.Range("E5").Value = -1.921949808943092
$ws.Range("F5").Value = 0.987
$ws.Range("H5").Value = 102
$ws.Range("I5").Value = 182
$ws.Range("J5").Value = 25.1
$ws.Range("K5").Value = 45.1
$ws.Range("L5").Value = 211.770852
$ws.Range("M5").Value = 173.890425
$ws.Range("N5").Value = -1.921949808943092
$ws.Range("O5").Value = -0.07817522212233576
$ws.Range("P5").Value = -1.843774586820756
$ws.Range("Q5").Value = -1.843774586820756
$ws.Range("T5").Value = 0.1418926829268293
$ws.Range("V5").Value = 0.0002529432437181515
$ws.Range("Z5").Value = -0.2616181228363033
$ws.Range("AB6").Value = -1476.076647556297
$ws.Range("AC6").Value = "umolO2/min/m2"
$ws.Range("AD6").Value = -1476.076647556297
$ws.Range("D6").Value = 293.5853077219285
$ws.Range("E6").Value = -1.967722579855463
$ws.Range("F6").Value = 0.98
$ws.Range("H6").Value = 101
$ws.Range("I6").Value = 181
$ws.Range("J6").Value = 24.91666666666667
$ws.Range("K6").Value = 44.91666666666666
$ws.Range("L6").Value = 243.146471
$ws.Range("M6").Value = 202.78998
$ws.Range("N6").Value = -1.967722579855463
$ws.Range("O6").Value = -0.07817522212233576
$ws.Range("P6").Value = -1.889547357733127
$ws.Range("Q6").Value = -1.889547357733127
$ws.Range("T6").Value = 0.1446439024390244
$ws.Range("V6").Value = 0.0001851607801792304
$ws.Range("Z6").Value = -0.2733115036658668
$ws.Range("AB7").Value = -764.489026552552
$ws.Range("AC7").Value = "umolO2/min/m2"
$ws.Range("AD7").Value = -764.489026552552
$ws.Range("D7").Value = 233.7285225591691
$ws.Range("E7").Value = -1.807070560072269
$ws.Range("F7").Value = 0.967
$ws.Range("H7").Value = 101
$ws.Range("I7").Value = 181
$ws.Range("J7").Value = 24.96666666666667
$ws.Range("K7").Value = 44.96666666666667
$ws.Range("L7").Value = 183.911286
$ws.Range("M7").Value = 154.640582
$ws.Range("N7").Value = -1.807070560072269
$ws.Range("O7").Value = -0.07817522212233576
$ws.Range("P7").Value = -1.728895337949933
$ws.Range("Q7").Value = -1.728895337949933
$ws.Range("T7").Value = 0.1429268292682927
$ws.Range("V7").Value = 0.0003232296608680373
$ws.Range("Z7").Value = -0.2471055287899174
$ws.Range("AB8").Value = "Inf"
$ws.Range("AC8").Value = "umolO2/min/m2"
$ws.Range("AD8").Value = "Inf"
$ws.Range("D8").Value = 185.4053589165027
$ws.Range("E8").Value = -0.07113783568202547
$ws.Range("F8").Value = 0.126
$ws.Range("H8").Value = 102
$ws.Range("I8").Value = 182
$ws.Range("J8").Value = 25.01666666666667
$ws.Range("K8").Value = 45.01666666666667
$ws.Range("L8").Value = 183.560071
$ws.Range("M8").Value = 182.511408
$ws.Range("N8").Value = -0.07113783568202547
$ws.Range("O8").Value = -0.07817522212233576
$ws.Range("P8").Value = 0.007037386440310281
$ws.Range("Q8").Value = 0.007037386440310281
$ws.Range("T8").Value = 0.1544
$ws.Range("V8").Value = 0
$ws.Range("Z8").Value = 0.001086572466383907
$ws.Range("AB9").Value = 1032.789634435439
$ws.Range("AC9").Value = "umolO2/min/m2"
$ws.Range("AD9").Value = 1032.789634435439
$ws.Range("D9").Value = 141.1527022944132
$ws.Range("E9").Value = 1.769327331454749
$ws.Range("F9").Value = 0.967
$ws.Range("I9").Value = 81
$ws.Range("K9").Value = 19.95
$ws.Range("M9").Value = 178.729241
$ws.Range("N9").Value = 1.769327331454749
$ws.Range("O9").Value = 0.02532265339143825
$ws.Range("P9").Value = 1.74400467806331
$ws.Range("Q9").Value = 1.74400467806331
$ws.Range("T9").Value = 0.1450048780487805
$ws.Range("V9").Value = 0.0002448603057459146
$ws.Range("Z9").Value = 0.252889185659073
$ws.Range("AB10").Value = 2114.00289476711
$ws.Range("AC10").Value = "umolO2/min/m2"
$ws.Range("AD10").Value = 2114.00289476711
$ws.Range("D10").Value = 157.1015984756967
$ws.Range("E10").Value = 2.134348235946302
$ws.Range("F10").Value = 0.982
$ws.Range("I10").Value = 81
$ws.Range("K10").Value = 19.98333333333333
$ws.Range("M10").Value = 204.186297
$ws.Range("N10").Value = 2.134348235946302
$ws.Range("O10").Value = 0.02532265339143825
$ws.Range("P10").Value = 2.109025582554863
$ws.Range("Q10").Value = 2.109025582554863
$ws.Range("T10").Value = 0.1492487804878049
$ws.Range("V10").Value = 0.0001488973818309612
$ws.Range("Z10").Value = 0.3147694962138956
$ws.Range("AB11").Value = 730.1297665824625
$ws.Range("AC11").Value = "umolO2/min/m2"
$ws.Range("AD11").Value = 730.1297665824625
$ws.Range("D11").Value = 129.9671140085762
$ws.Range("E11").Value = 1.129912095768344
$ws.Range("F11").Value = 0.629
$ws.Range("I11").Value = 82
$ws.Range("K11").Value = 20.05
$ws.Range("M11").Value = 151.134486
$ws.Range("N11").Value = 1.129912095768344
$ws.Range("O11").Value = 0.02532265339143825
$ws.Range("P11").Value = 1.104589442376905
$ws.Range("Q11").Value = 1.104589442376905
$ws.Range("T11").Value = 0.1469268292682927
$ws.Range("V11").Value = 0.0002222807942365138
$ws.Range("Z11").Value = 0.1622938244116702
$ws.Range("AB12").Value = 1182.73577054592
$ws.Range("AC12").Value = "umolO2/min/m2"
$ws.Range("AD12").Value = 1182.73577054592
$ws.Range("D12").Value = 170.6055724988584
$ws.Range("E12").Value = 2.133712008589159
$ws.Range("F12").Value = 0.982
$ws.Range("I12").Value = 82
$ws.Range("K12").Value = 20.1
$ws.Range("M12").Value = 215.160337
$ws.Range("N12").Value = 2.133712008589159
$ws.Range("O12").Value = 0.02532265339143825
$ws.Range("P12").Value = 2.108389355197721
$ws.Range("Q12").Value = 2.108389355197721
$ws.Range("T12").Value = 0.1418926829268293
$ws.Range("V12").Value = 0.0002529432437181515
$ws.Range("Z12").Value = 0.2991650222633723
$ws.Range("AB13").Value = 1977.806163768549
$ws.Range("AC13").Value = "umolO2/min/m2"
$ws.Range("AD13").Value = 1977.806163768549
$ws.Range("D13").Value = 193.543223554161
$ws.Range("E13").Value = 2.557141320832704
$ws.Range("F13").Value = 0.965
$ws.Range("I13").Value = 81
$ws.Range("K13").Value = 19.91666666666667
$ws.Range("M13").Value = 247.69543
$ws.Range("N13").Value = 2.557141320832704
$ws.Range("O13").Value = 0.02532265339143825
$ws.Range("P13").Value = 2.531818667441266
$ws.Range("Q13").Value = 2.531818667441266
$ws.Range("T13").Value = 0.1446439024390244
$ws.Range("V13").Value = 0.0001851607801792304
$ws.Range("Z13").Value = 0.3662121323266753
$ws.Range("AB14").Value = 808.1247811426513
$ws.Range("AC14").Value = "umolO2/min/m2"
$ws.Range("AD14").Value = 808.1247811426513
$ws.Range("D14").Value = 149.1277181191609
$ws.Range("E14").Value = 1.852900444664414
$ws.Range("F14").Value = 0.969
$ws.Range("I14").Value = 81
$ws.Range("K14").Value = 19.96666666666667
$ws.Range("M14").Value = 185.147121
$ws.Range("N14").Value = 1.852900444664414
$ws.Range("O14").Value = 0.02532265339143825
$ws.Range("P14").Value = 1.827577791272976
$ws.Range("Q14").Value = 1.827577791272976
$ws.Range("T14").Value = 0.1429268292682927
$ws.Range("V14").Value = 0.0003232296608680373
$ws.Range("Z14").Value = 0.2612098989477961
$ws.Range("AB15").Value = "Inf"
$ws.Range("AC15").Value = "umolO2/min/m2"
$ws.Range("AD15").Value = "Inf"
$ws.Range("D15").Value = 183.3523673100485
$ws.Range("E15").Value = 0.02854873738546629
$ws.Range("F15").Value = 0.0113
$ws.Range("I15").Value = 82
$ws.Range("K15").Value = 20.01666666666667
$ws.Range("M15").Value = 184.792976
$ws.Range("N15").Value = 0.02854873738546629
$ws.Range("O15").Value = 0.02532265339143825
$ws.Range("P15").Value = 0.003226083994028035
$ws.Range("Q15").Value = 0.003226083994028035
$ws.Range("T15").Value = 0.1544
$ws.Range("V15").Value = 0
$ws.Range("Z15").Value = 0.0004981073686779286
